$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.413144
$ws.Range("H2").Value = 1.239432
$ws.Range("I2").Value = 0.4553782032534783
$ws.Range("J2").Value = 0.4553782032534782
$ws.Range("Q2").Value = 0.04118894193866666
$ws.Range("R2").Value = 0.370700477448
$ws.Range("S2").Value = 0.4553782032534783
$ws.Range("T2").Value = 0.4553782032534782

# Row 3 updates
$ws.Range("I3").Value = 0.3895918235379703
$ws.Range("J3").Value = 0.3895918235379702
$ws.Range("R3").Value = 0.3171470965530001
$ws.Range("S3").Value = 0.3895918235379703
$ws.Range("T3").Value = 0.3895918235379702

# Row 4 updates
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1406516666666666
$ws.Range("H4").Value = 0.421955
$ws.Range("I4").Value = 0.1550299732085515
$ws.Range("J4").Value = 0.1550299732085515
$ws.Range("Q4").Value = 0.01402245544388889
$ws.Range("R4").Value = 0.126202098995
$ws.Range("S4").Value = 0.1550299732085515
$ws.Range("T4").Value = 0.1550299732085515
